$wb = $excel.ActiveWorkbook

# --- Sheet 1: summ55077076 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ55077076"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = -2162.979324538224
$ws.Cells.Item(2,3).Value = 0.5832351171227304
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -1956.57744153977
$ws.Cells.Item(3,3).Value = 0.453089369859632
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 447.8059771910899
$ws.Cells.Item(4,3).Value = 0.7014024246142423
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -977.6163418770144
$ws.Cells.Item(5,3).Value = 0.0440174230489936
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 528.7599739931343
$ws.Cells.Item(6,3).Value = 0.3761719637963514
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = -174.7004586806873
$ws.Cells.Item(7,3).Value = 0.7767904245226606
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = -54.17121621971927
$ws.Cells.Item(8,3).Value = 0.9307542619528214
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 316.4715990932069
$ws.Cells.Item(9,3).Value = 0.1087488647169929
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1870.295981595641
$ws.Cells.Item(10,3).Value = 0.00001039392017342865
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -33.48475522836647
$ws.Cells.Item(11,3).Value = 0.09530322395215426
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 249.8812724599509
$ws.Cells.Item(12,3).Value = 0.238301709462112
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 882.5302702796218
$ws.Cells.Item(13,3).Value = 0.00001075815875222905
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = -0.06509053596760833
$ws.Cells.Item(14,3).Value = 0.7750003448775094
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.0001421345314552166
$ws.Cells.Item(15,3).Value = 0.2183714720688381
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 131.3843915843905
$ws.Cells.Item(16,3).Value = 0.004614023779743567
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 54.79714624330452
$ws.Cells.Item(17,3).Value = 0.02391706815553074
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -9796.591316089329
$ws.Cells.Item(18,3).Value = 0.006817555895811732
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -7253.903041275493
$ws.Cells.Item(19,3).Value = 0.1563244770228284
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -2947.551234406144
$ws.Cells.Item(20,3).Value = 0.2092396367128242

# --- Sheet 2: summ55267027 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ55267027"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = -1260.970791802835
$ws.Cells.Item(2,3).Value = 0.7578865221717519
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -2919.47328752046
$ws.Cells.Item(3,3).Value = 0.1969402422814156
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 283.8702675671981
$ws.Cells.Item(4,3).Value = 0.8046733532174313
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -858.9298786586809
$ws.Cells.Item(5,3).Value = 0.08259939220076265
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 467.5281762131926
$ws.Cells.Item(6,3).Value = 0.4397520557340604
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = -343.9925670138518
$ws.Cells.Item(7,3).Value = 0.5834767751300796
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = -88.72958042163657
$ws.Cells.Item(8,3).Value = 0.8901041288300461
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 343.7933932933638
$ws.Cells.Item(9,3).Value = 0.0841433477841682
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1795.419211369705
$ws.Cells.Item(10,3).Value = 0.00003050236359605737
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -31.35565098387175
$ws.Cells.Item(11,3).Value = 0.1291418852034751
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 337.7536726697313
$ws.Cells.Item(12,3).Value = 0.122891310134152
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 869.4558904419239
$ws.Cells.Item(13,3).Value = 0.00002578422094954419
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = -0.05942168754003412
$ws.Cells.Item(14,3).Value = 0.7970601788892798
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.0001664371736413507
$ws.Cells.Item(15,3).Value = 0.1568550324286831
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 102.1695814069161
$ws.Cells.Item(16,3).Value = 0.02985446380823706
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 48.22797437376347
$ws.Cells.Item(17,3).Value = 0.05278723838560756
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -10418.08038189263
$ws.Cells.Item(18,3).Value = 0.003702246166738263
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -5760.74193867565
$ws.Cells.Item(19,3).Value = 0.2650281999081638
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -917.6322106572363
$ws.Cells.Item(20,3).Value = 0.6965875324640458

# --- Sheet 3: summ55472550 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ55472550"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = 338.0476323535599
$ws.Cells.Item(2,3).Value = 0.9331840001284362
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -1661.627186258003
$ws.Cells.Item(3,3).Value = 0.4875297234936981
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 1049.735055884742
$ws.Cells.Item(4,3).Value = 0.3757994117451073
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -887.8651889562018
$ws.Cells.Item(5,3).Value = 0.06556297223801394
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 575.9696846985828
$ws.Cells.Item(6,3).Value = 0.3373093737542314
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = -180.7007310276472
$ws.Cells.Item(7,3).Value = 0.7699381991882942
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = 160.2021921417289
$ws.Cells.Item(8,3).Value = 0.7996103821108979
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 385.1525030534422
$ws.Cells.Item(9,3).Value = 0.05507464629850244
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1845.660266811134
$ws.Cells.Item(10,3).Value = 0.00001456812058989902
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -37.75700774187823
$ws.Cells.Item(11,3).Value = 0.06459685307064947
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 394.6753590216238
$ws.Cells.Item(12,3).Value = 0.06038437407732916
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 859.9064865224678
$ws.Cells.Item(13,3).Value = 0.00001972927658032867
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = -0.1888809707851697
$ws.Cells.Item(14,3).Value = 0.4023271623628204
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.0001780380192987479
$ws.Cells.Item(15,3).Value = 0.122617099414513
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 98.42698701800077
$ws.Cells.Item(16,3).Value = 0.03416718924394788
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 34.42381352673195
$ws.Cells.Item(17,3).Value = 0.1600255009195364
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -10722.74185215208
$ws.Cells.Item(18,3).Value = 0.002898731584378516
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -4050.921476550283
$ws.Cells.Item(19,3).Value = 0.4219660120723648
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -980.5849755290674
$ws.Cells.Item(20,3).Value = 0.6721630751262221

# --- Sheet 4: summ55667129 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ55667129"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = 446.7602255796219
$ws.Cells.Item(2,3).Value = 0.9129210949848102
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -2925.378285784338
$ws.Cells.Item(3,3).Value = 0.1971583492133837
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 1123.816604879433
$ws.Cells.Item(4,3).Value = 0.3429840553592579
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -1048.715443025896
$ws.Cells.Item(5,3).Value = 0.03414740223483834
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 607.9377788870365
$ws.Cells.Item(6,3).Value = 0.3167704086573179
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = -237.7084271848356
$ws.Cells.Item(7,3).Value = 0.7065612598805343
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = -281.3715823464201
$ws.Cells.Item(8,3).Value = 0.6596277427467263
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 326.2851501963627
$ws.Cells.Item(9,3).Value = 0.1052731450846236
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1852.262754644806
$ws.Cells.Item(10,3).Value = 0.00001813941544349227
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -33.78161651082612
$ws.Cells.Item(11,3).Value = 0.1021698547233887
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 195.5148091386729
$ws.Cells.Item(12,3).Value = 0.3643809215207743
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 889.3809184776756
$ws.Cells.Item(13,3).Value = 0.0000143809127544995
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = -0.1131720591381013
$ws.Cells.Item(14,3).Value = 0.6255527457715105
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.0001567631589791746
$ws.Cells.Item(15,3).Value = 0.1823063234515749
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 114.4487706693059
$ws.Cells.Item(16,3).Value = 0.01564299900255494
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 40.67435139881236
$ws.Cells.Item(17,3).Value = 0.1001061003575959
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -11319.2153475951
$ws.Cells.Item(18,3).Value = 0.001803648051710673
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -6416.302092866699
$ws.Cells.Item(19,3).Value = 0.2157036036285415
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -2205.126977085284
$ws.Cells.Item(20,3).Value = 0.3550883792032816

# --- Sheet 5: summ55861763 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ55861763"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = -1472.83694328345
$ws.Cells.Item(2,3).Value = 0.7127400185378311
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -2701.260317082908
$ws.Cells.Item(3,3).Value = 0.2239937929687487
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 973.5421110937771
$ws.Cells.Item(4,3).Value = 0.3875378803425975
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -1172.683272168054
$ws.Cells.Item(5,3).Value = 0.01601320828352553
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 430.8713819013109
$ws.Cells.Item(6,3).Value = 0.4700561182445833
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = 22.51985573657907
$ws.Cells.Item(7,3).Value = 0.9710777085599394
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = -154.8262888645804
$ws.Cells.Item(8,3).Value = 0.8064557634261567
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 467.7785345254968
$ws.Cells.Item(9,3).Value = 0.01521674929851043
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1898.844851343635
$ws.Cells.Item(10,3).Value = 0.00000721069953527011
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -19.9479243488081
$ws.Cells.Item(11,3).Value = 0.3185570922315045
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 265.5843901945518
$ws.Cells.Item(12,3).Value = 0.2067546007456944
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 825.3651661579548
$ws.Cells.Item(13,3).Value = 0.00004093418386019214
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = -0.08944179263883466
$ws.Cells.Item(14,3).Value = 0.7041775893526314
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.00013285676267865
$ws.Cells.Item(15,3).Value = 0.2445222685619411
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 98.6950545683622
$ws.Cells.Item(16,3).Value = 0.03335515186487176
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 49.52954393090396
$ws.Cells.Item(17,3).Value = 0.04224668801438331
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -9435.34213194832
$ws.Cells.Item(18,3).Value = 0.006914268626912614
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -5733.060444978331
$ws.Cells.Item(19,3).Value = 0.2681060551240345
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -1433.414604463182
$ws.Cells.Item(20,3).Value = 0.5368649332584043

# --- Sheet 6: summ56050671 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ56050671"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = 182.2734595813145
$ws.Cells.Item(2,3).Value = 0.961829049971199
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -6110.482078147326
$ws.Cells.Item(3,3).Value = 0.009330363662335759
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 1512.821032215799
$ws.Cells.Item(4,3).Value = 0.1714067295299627
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -1033.394709581032
$ws.Cells.Item(5,3).Value = 0.0290865244279693
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 135.3146736472374
$ws.Cells.Item(6,3).Value = 0.8156517426014345
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = -373.6506417481567
$ws.Cells.Item(7,3).Value = 0.5354949552188843
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = -195.0080573670322
$ws.Cells.Item(8,3).Value = 0.7521867433063866
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 442.5216195920843
$ws.Cells.Item(9,3).Value = 0.01945345176820592
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1588.740136133896
$ws.Cells.Item(10,3).Value = 0.0001102556729981034
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -28.60248786680152
$ws.Cells.Item(11,3).Value = 0.1452626009434535
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 170.0795492665736
$ws.Cells.Item(12,3).Value = 0.4186018618151814
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 796.7210449596555
$ws.Cells.Item(13,3).Value = 0.00004893649200414999
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = 0.09570885438653125
$ws.Cells.Item(14,3).Value = 0.6634464617519287
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.0001541844518303371
$ws.Cells.Item(15,3).Value = 0.1771589224383964
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 104.1533683216869
$ws.Cells.Item(16,3).Value = 0.02083758082226005
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 40.76288287419138
$ws.Cells.Item(17,3).Value = 0.08310565623859333
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -9353.265819116445
$ws.Cells.Item(18,3).Value = 0.00798068806684665
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -8254.075032635737
$ws.Cells.Item(19,3).Value = 0.09626290472430303
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -3641.023672379536
$ws.Cells.Item(20,3).Value = 0.1066859505206911

# --- Sheet 7: summ56246027 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ56246027"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = -379.9501930656552
$ws.Cells.Item(2,3).Value = 0.9243970289971795
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -2299.908596342279
$ws.Cells.Item(3,3).Value = 0.335952436350619
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 704.0526861499411
$ws.Cells.Item(4,3).Value = 0.5503403798379722
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -1056.06768259099
$ws.Cells.Item(5,3).Value = 0.03083191781340099
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 990.5966706240624
$ws.Cells.Item(6,3).Value = 0.09681139180336891
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = 170.297343807739
$ws.Cells.Item(7,3).Value = 0.7827534596613667
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = 31.54853428528997
$ws.Cells.Item(8,3).Value = 0.9602541803882135
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 271.4549894839393
$ws.Cells.Item(9,3).Value = 0.165485141042949
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1972.738878236359
$ws.Cells.Item(10,3).Value = 0.000003394203426597894
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -33.09500079676697
$ws.Cells.Item(11,3).Value = 0.1036311769842848
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 355.1139053481292
$ws.Cells.Item(12,3).Value = 0.09593810378421853
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 873.7883624090395
$ws.Cells.Item(13,3).Value = 0.00001448597520772066
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = -0.03700922529318178
$ws.Cells.Item(14,3).Value = 0.8719927958277485
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.0002107351928149421
$ws.Cells.Item(15,3).Value = 0.07289194887021871
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 101.1119540010184
$ws.Cells.Item(16,3).Value = 0.02900820681976467
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 44.03321386894584
$ws.Cells.Item(17,3).Value = 0.07072758069435156
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -11167.12276861967
$ws.Cells.Item(18,3).Value = 0.002236863042656039
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -7411.092895176283
$ws.Cells.Item(19,3).Value = 0.1461006082139285
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -881.1212865340657
$ws.Cells.Item(20,3).Value = 0.7068217106044358

# --- Sheet 8: summ56435032 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ56435032"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = 186.4449310554205
$ws.Cells.Item(2,3).Value = 0.9624860226042778
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -3101.916740910754
$ws.Cells.Item(3,3).Value = 0.2394368855191219
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 734.7209257839014
$ws.Cells.Item(4,3).Value = 0.5284804401094831
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -976.0758816466978
$ws.Cells.Item(5,3).Value = 0.04523830403550039
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 615.2859469896409
$ws.Cells.Item(6,3).Value = 0.2995921846752328
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = 56.6265518753747
$ws.Cells.Item(7,3).Value = 0.9271777049229086
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = 142.3639286435853
$ws.Cells.Item(8,3).Value = 0.8219378393322557
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 359.2267306985063
$ws.Cells.Item(9,3).Value = 0.06405055509682658
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1910.948046093477
$ws.Cells.Item(10,3).Value = 0.000007569215820940635
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -29.37638230368781
$ws.Cells.Item(11,3).Value = 0.1428212160637385
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 394.2199219410722
$ws.Cells.Item(12,3).Value = 0.06205702033775525
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 807.8528483748307
$ws.Cells.Item(13,3).Value = 0.00005300481584907108
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = -0.00512054273780238
$ws.Cells.Item(14,3).Value = 0.9820012193565242
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.0001851638245639747
$ws.Cells.Item(15,3).Value = 0.1139008980452475
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 96.3986205118093
$ws.Cells.Item(16,3).Value = 0.03749977093788528
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 37.58203427892156
$ws.Cells.Item(17,3).Value = 0.1176055803456772
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -10321.31638959795
$ws.Cells.Item(18,3).Value = 0.00330938462587107
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -7410.877632199993
$ws.Cells.Item(19,3).Value = 0.1400769690262201
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -606.4626745629193
$ws.Cells.Item(20,3).Value = 0.7963882049012662

# --- Sheet 9: summ56625575 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ56625575"
$ws.Cells.Item(2,1).Value = "Intercept"
$ws.Cells.Item(2,2).Value = 1983.059147996308
$ws.Cells.Item(2,3).Value = 0.610056037933211
$ws.Cells.Item(3,1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3,2).Value = -3352.93496312556
$ws.Cells.Item(3,3).Value = 0.1261374523799426
$ws.Cells.Item(4,1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4,2).Value = 1140.537935993043
$ws.Cells.Item(4,3).Value = 0.3115296148003179
$ws.Cells.Item(5,1).Value = "Education[T.University]"
$ws.Cells.Item(5,2).Value = -1293.666087759981
$ws.Cells.Item(5,3).Value = 0.00709258227763781
$ws.Cells.Item(6,1).Value = "Season[T.Spring]"
$ws.Cells.Item(6,2).Value = 85.62713215273567
$ws.Cells.Item(6,3).Value = 0.8844493315271553
$ws.Cells.Item(7,1).Value = "Season[T.Summer]"
$ws.Cells.Item(7,2).Value = -441.42313001613
$ws.Cells.Item(7,3).Value = 0.4687396117306635
$ws.Cells.Item(8,1).Value = "Season[T.Winter]"
$ws.Cells.Item(8,2).Value = -335.2190470801816
$ws.Cells.Item(8,3).Value = 0.5893444924842917
$ws.Cells.Item(9,1).Value = "HHSize"
$ws.Cells.Item(9,2).Value = 522.9600146536397
$ws.Cells.Item(9,3).Value = 0.007991983722532355
$ws.Cells.Item(10,1).Value = "Sex"
$ws.Cells.Item(10,2).Value = -1485.783541707172
$ws.Cells.Item(10,3).Value = 0.0003603440613907803
$ws.Cells.Item(11,1).Value = "Age"
$ws.Cells.Item(11,2).Value = -42.31802489502659
$ws.Cells.Item(11,3).Value = 0.03368788638071963
$ws.Cells.Item(12,1).Value = "DistSubcenter_res"
$ws.Cells.Item(12,2).Value = 364.3037457004467
$ws.Cells.Item(12,3).Value = 0.0816640241474337
$ws.Cells.Item(13,1).Value = "DistCenter_res"
$ws.Cells.Item(13,2).Value = 887.847898644101
$ws.Cells.Item(13,3).Value = 0.000009651348705932376
$ws.Cells.Item(14,1).Value = "UrbPopDensity_res"
$ws.Cells.Item(14,2).Value = -0.03306932549440245
$ws.Cells.Item(14,3).Value = 0.8845400008254145
$ws.Cells.Item(15,1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(15,2).Value = 0.0001917809965348715
$ws.Cells.Item(15,3).Value = 0.09339257909716771
$ws.Cells.Item(16,1).Value = "IntersecDensity_res"
$ws.Cells.Item(16,2).Value = 90.42351068621703
$ws.Cells.Item(16,3).Value = 0.04295274120280176
$ws.Cells.Item(17,1).Value = "street_length_res"
$ws.Cells.Item(17,2).Value = 23.73005756026974
$ws.Cells.Item(17,3).Value = 0.3159785368678439
$ws.Cells.Item(18,1).Value = "LU_Comm_res"
$ws.Cells.Item(18,2).Value = -11735.72571669738
$ws.Cells.Item(18,3).Value = 0.0007253149354887275
$ws.Cells.Item(19,1).Value = "LU_UrbFab_res"
$ws.Cells.Item(19,2).Value = -5532.987990333781
$ws.Cells.Item(19,3).Value = 0.2680086523546668
$ws.Cells.Item(20,1).Value = "bike_lane_share_res"
$ws.Cells.Item(20,2).Value = -1094.905194406664
$ws.Cells.Item(20,3).Value = 0.6350533597480745
